$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (r) for columns B..G (ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(-0.2626095943911887, 1.949383735764183, 16.52272044676487, 4.064814933888734, 4.147487754911329, 23)
    3  = @(-0.01304390213147815, 1.772198340847518, 10.9418263788771, 3.307843161166669, 3.385659091639844, 22)
    4  = @(-0.660492300833325, 1.197293372571959, 5.006510981504795, 2.237523403565825, 2.190610405148969, 21)
    5  = @(-0.05106328621311902, 0.7433681293355469, 1.840802665129192, 1.356761830657537, 1.391022042252352, 20)
    6  = @(-0.06395863908153489, 0.8000673122416484, 1.420101517000478, 1.191680123607203, 1.222570274741316, 19)
    7  = @(-0.03738556241501667, 0.71332745744033, 0.9386152008072798, 0.9688215526129049, 0.9961666341101381, 18)
    8  = @(0.06091063799584202, 0.6343285541743425, 0.7182370180715076, 0.8474886536535505, 0.8713121394001694, 17)
    9  = @(0.1648078376396195, 0.5511977509794913, 0.4449526034813988, 0.6670476770676882, 0.6675655114989703, 16)
    10 = @(0.1934705323399399, 0.5837592819091525, 0.5633658838424578, 0.7505770339162116, 0.7506675864274018, 15)
    11 = @(0.2495182209949923, 0.5857948382811914, 0.4537182255294631, 0.6735860936283224, 0.6492850719707854, 14)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
